$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Summary": update aggregate stats to reflect the new trade
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1500.1    # Current Capital
$wsSummary.Range("B4").Value = 0.1       # Total P&L $
$wsSummary.Range("B5").Value = 0.5       # Total P&L %
$wsSummary.Range("B6").Value = 4         # Total Trades
$wsSummary.Range("B7").Value = 2         # Winning Trades
$wsSummary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 6)
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 100.1      # Capital
$wsStatus.Range("D6").Value = 4          # Trades
$wsStatus.Range("E6").Value = 0.1        # P&L $
$wsStatus.Range("F6").Value = 0.1        # P&L %
$wsStatus.Range("G6").Value = 50         # Win Rate %

# ---------------------------------------------------------------
# Sheet "All Trades": append new trade #4 as row 5
# ---------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("A5").Value = 4
# Force Date/Time columns to be written as literal text (matching rows 2-4)
# instead of being auto-parsed into Excel date/time serial numbers.
$wsAllTrades.Range("B5").NumberFormat = "@"
$wsAllTrades.Range("B5").Value = "2026-02-17"
$wsAllTrades.Range("B5").Style = "Normal"
$wsAllTrades.Range("C5").NumberFormat = "@"
$wsAllTrades.Range("C5").Value = "23:52:09"
$wsAllTrades.Range("C5").Style = "Normal"
$wsAllTrades.Range("D5").Value = "MarketMaking"
$wsAllTrades.Range("E5").Value = "UP"
$wsAllTrades.Range("F5").Value = 0.9
$wsAllTrades.Range("G5").Value = 0.91
$wsAllTrades.Range("H5").Value = "CLOSED"
$wsAllTrades.Range("I5").Value = 1.1111
$wsAllTrades.Range("J5").Value = 0.01
$wsAllTrades.Range("K5").Value = 100.1
$wsAllTrades.Range("L5").Value = 0
$wsAllTrades.Range("M5").Value = 0
$wsAllTrades.Range("N5").Value = 0.6
$wsAllTrades.Range("O5").Value = "Normal spread capture: 19600 bps"
$wsAllTrades.Range("P5").Value = "early_exit"
$wsAllTrades.Range("Q5").Value = 0.14

# ---------------------------------------------------------------
# Sheet "MarketMaking": append the same new trade #4 as row 5
# ---------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A5").Value = 4
# Force Date/Time columns to be written as literal text (matching rows 2-4)
# instead of being auto-parsed into Excel date/time serial numbers.
$wsMM.Range("B5").NumberFormat = "@"
$wsMM.Range("B5").Value = "2026-02-17"
$wsMM.Range("B5").Style = "Normal"
$wsMM.Range("C5").NumberFormat = "@"
$wsMM.Range("C5").Value = "23:52:09"
$wsMM.Range("C5").Style = "Normal"
$wsMM.Range("D5").Value = "MarketMaking"
$wsMM.Range("E5").Value = "UP"
$wsMM.Range("F5").Value = 0.9
$wsMM.Range("G5").Value = 0.91
$wsMM.Range("H5").Value = "CLOSED"
$wsMM.Range("I5").Value = 1.1111
$wsMM.Range("J5").Value = 0.01
$wsMM.Range("K5").Value = 100.1
$wsMM.Range("L5").Value = 0
$wsMM.Range("M5").Value = 0
$wsMM.Range("N5").Value = 0.6
$wsMM.Range("O5").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("P5").Value = "early_exit"
$wsMM.Range("Q5").Value = 0.14
